# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets
# to match the freshly scraped snapshot (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 15425
$ws.Range("F6").Value  = 417
$ws.Range("F9").Value  = 15361
$ws.Range("F10").Value = 48
$ws.Range("F11").Value = 8955
$ws.Range("F12").Value = 365
$ws.Range("F14").Value = 1009
$ws.Range("F15").Value = 81
$ws.Range("F16").Value = 194
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 17
$ws.Range("F21").Value = 538
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 74
$ws.Range("F32").Value = 52
$ws.Range("F33").Value = 37
$ws.Range("F35").Value = 301
$ws.Range("F36").Value = 440
$ws.Range("F38").Value = 5477

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 15425
$ws.Range("F6").Value  = 417
$ws.Range("F9").Value  = 15361
$ws.Range("F10").Value = 48
$ws.Range("F11").Value = 8955
$ws.Range("F12").Value = 365
$ws.Range("F14").Value = 1009
$ws.Range("F15").Value = 81
$ws.Range("F16").Value = 194
$ws.Range("F18").Value = 193
$ws.Range("F19").Value = 17
$ws.Range("F21").Value = 538
$ws.Range("F27").Value = 21
$ws.Range("F28").Value = 74
$ws.Range("F34").Value = 52
$ws.Range("F35").Value = 37
$ws.Range("F37").Value = 301
$ws.Range("F38").Value = 440
$ws.Range("F40").Value = 5477
